$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.732.10'
$ws.Range('E2').Value = '  -2.77%  '
$ws.Range('D3').Value = '2.903.01'
$ws.Range('E3').Value = '  -3.88%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -2.81%  '
$ws.Range('D9').Value = '2.901.09'
$ws.Range('E9').Value = '  -3.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.33%  '
$ws.Range('E11').Value = '  -4.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.448'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.70%  '
$ws.Range('E13').Value = '  -3.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.97'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('E15').Value = '  +0.70%  '
$ws.Range('D16').Value = '3.382.53'
$ws.Range('E16').Value = '  -3.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.80'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.26%  '
$ws.Range('D18').Value = '60.598.95'
$ws.Range('E18').Value = '  -2.75%  '
$ws.Range('D19').Value = '2.903.18'
$ws.Range('E19').Value = '  -3.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '427.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.36%  '
$ws.Range('E22').Value = '  -3.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.63'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.08'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.21'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.22%  '
$ws.Range('E32').Value = '  -3.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.47'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.11%  '
$ws.Range('E34').Value = '  -3.34%  '
$ws.Range('D35').Value = '0.0₃0835'
$ws.Range('E35').Value = '  -2.08%  '
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.66'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.39%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.03'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.28'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  -2.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.122'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.72'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.49%  '
$ws.Range('E43').Value = '  +2.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0346'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '371.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '133.78'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('D48').Value = '2.656.44'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.45%  '
$ws.Range('E51').Value = '  -1.47%  '
